$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.652.58"
$ws.Range("E2").Value = "  +3.65%  "
$ws.Range("D3").Value = "3.344.96"
$ws.Range("E3").Value = "  +8.55%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.23"
$ws.Range("E5").Value = "  +10.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.04"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.13"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.342.88"
$ws.Range("E10").Value = "  +8.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.798"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("D13").Value = "97.394.65"
$ws.Range("E13").Value = "  +3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.60"
$ws.Range("E14").Value = "  +5.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000247"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "3.962.62"
$ws.Range("E16").Value = "  +8.49%  "
$ws.Range("E17").Value = "  +5.03%  "
$ws.Range("D18").Value = "3.344.02"
$ws.Range("E18").Value = "  +8.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.62"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.97"
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.89"
$ws.Range("E21").Value = "  +10.00%  "
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000208"
$ws.Range("E22").Value = "  +9.70%  "
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.20"
$ws.Range("E24").Value = "  +4.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.63"
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.23"
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("D28").Value = "3.525.59"
$ws.Range("E28").Value = "  +8.39%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.38"
$ws.Range("E35").Value = "  +7.30%  "
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "513.05"
$ws.Range("E38").Value = "  +10.41%  "
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.82"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.449"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.32"
$ws.Range("E43").Value = "  +7.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.53"
$ws.Range("E44").Value = "  -5.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.788"
$ws.Range("E45").Value = "  +16.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.03"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("E48").Value = "  +4.22%  "
$ws.Range("E49").Value = "  +7.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.47"
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.804"
$ws.Range("E51").Value = "  +11.71%  "
